$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "03cec64e-63f2-41db-ab9e-c54c60418425"
$ws.Range("B6").Value = "PVFCKC"
$ws.Range("C6").Value = 1746862626.118984
$ws.Range("D6").Value = $false
